$wb = $excel.ActiveWorkbook

# Header row is identical across all 15 sheets; same for every worksheet.
$headers = @{
  'A1' = 'Ratings'
  'B1' = 'API_and_Integration_Support'
  'C1' = 'Pricing_Details'
  'D1' = 'Deployment_Support'
  'E1' = 'Customer_Support_Options'
  'F1' = 'Training_Platforms'
  'G1' = 'Vendor_Details'
  'H1' = 'Features'
}

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
  $ws = $wb.Worksheets.Item($i)
  foreach ($cell in $headers.Keys) {
    $ws.Range($cell).Value = $headers[$cell]
  }
}

# sheet1.xml
$ws = $wb.Worksheets.Item(1)
$ws.Range('A2').Value = '{"Total_reviews":2500,"Ease_of_use":4.5,"Features":4.7,"Design":4.2,"Support":4}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":true,"Cozyroc SSIS+ Suite":false,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other_available_integrations":"Hundreds of integrations available through the ServiceNow store"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on usage and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":true,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In-person_training":true}'
$ws.Range('G2').Value = '{"Company_name":"ServiceNow","Year_founded":2004,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","IT Asset Management","Workflow Automation","Reporting and Analytics","Cloud Management","Security Management","Customer Service Management"]'

# sheet2.xml
$ws = $wb.Worksheets.Item(2)
$ws.Range('A2').Value = '{"Total_reviews":1500,"Ease_of_use":4.2,"Features":4,"Design":3.8,"Support":3.7}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In-person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"SolarWinds","Year_founded":1999,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","Workflow Automation","Reporting and Analytics"]'

# sheet3.xml
$ws = $wb.Worksheets.Item(3)
$ws.Range('A2').Value = '{"Total_reviews":1200,"Ease_of_use":4,"Features":4.2,"Design":3.5,"Support":3.8}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In-person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"ManageEngine","Year_founded":2001,"Country":"India"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics"]'

# sheet4.xml
$ws = $wb.Worksheets.Item(4)
$ws.Range('A2').Value = '{"Total_reviews":800,"Ease_of_use":4.3,"Features":4.1,"Design":3.9,"Support":4}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In-person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"TOPdesk","Year_founded":1993,"Country":"Netherlands"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics"]'

# sheet5.xml
$ws = $wb.Worksheets.Item(5)
$ws.Range('A2').Value = '{"Total_reviews":500,"Ease_of_use":4,"Features":4.3,"Design":3.8,"Support":3.7}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In-person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"SymphonyAI","Year_founded":2017,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","AI-powered automation","Predictive analytics"]'

# sheet6.xml
$ws = $wb.Worksheets.Item(6)
$ws.Range('A2').Value = '{"Total_reviews":2000,"Ease_of_use":4.2,"Features":4.4,"Design":3.9,"Support":3.8}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":true,"Cozyroc SSIS+ Suite":false,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other_available_integrations":"Hundreds of integrations available through the Atlassian Marketplace"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":true}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In-person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"Atlassian","Year_founded":2002,"Country":"Australia"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","Workflow Automation","Reporting and Analytics","Project Management Integration","Agile Development Integration"]'

# sheet7.xml
$ws = $wb.Worksheets.Item(7)
$ws.Range('A2').Value = '{"Total_reviews":700,"Ease_of_use":4.1,"Features":4.3,"Design":3.7,"Support":3.9}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"Cherwell Software","Year_founded":1997,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics"]'

# sheet8.xml
$ws = $wb.Worksheets.Item(8)
$ws.Range('A2').Value = '{"Total_reviews":1800,"Ease_of_use":4.5,"Features":4.3,"Design":4.1,"Support":4.2}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":true,"Assess360":false,"BigID":true,"Cozyroc SSIS+ Suite":false,"CloudHub":true,"Elastic Observability":true,"Exalate":true,"Incydr":true,"Nexpose":true,"Other_available_integrations":"Hundreds of integrations available through the Freshworks Marketplace"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":true,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"Freshworks","Year_founded":2010,"Country":"India"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Customer Service Management","AI-powered automation"]'

# sheet9.xml
$ws = $wb.Worksheets.Item(9)
$ws.Range('A2').Value = '{"Total_reviews":1000,"Ease_of_use":4.2,"Features":4.1,"Design":3.8,"Support":4}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":true,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":true,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"SysAid Technologies","Year_founded":1999,"Country":"Israel"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Remote Desktop Support","Self-service Portal"]'

# sheet10.xml
$ws = $wb.Worksheets.Item(10)
$ws.Range('A2').Value = '{"Total_reviews":600,"Ease_of_use":3.8,"Features":4.2,"Design":3.5,"Support":3.7}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"BMC Software","Year_founded":1980,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Service Level Management"]'

# sheet11.xml
$ws = $wb.Worksheets.Item(11)
$ws.Range('A2').Value = '{"Total_reviews":600,"Ease_of_use":4,"Features":4.2,"Design":3.9,"Support":3.8}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"Ivanti","Year_founded":1994,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Endpoint Management","Unified Endpoint Management"]'

# sheet12.xml
$ws = $wb.Worksheets.Item(12)
$ws.Range('A2').Value = '{"Total_reviews":400,"Ease_of_use":4,"Features":4.1,"Design":3.7,"Support":3.9}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"EV Technologies","Year_founded":1996,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Service Level Management","ITIL Process Management"]'

# sheet13.xml
$ws = $wb.Worksheets.Item(13)
$ws.Range('A2').Value = '{"Total_reviews":900,"Ease_of_use":4.1,"Features":3.9,"Design":3.7,"Support":3.6}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"SolarWinds","Year_founded":1999,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","Workflow Automation","Reporting and Analytics"]'

# sheet14.xml
$ws = $wb.Worksheets.Item(14)
$ws.Range('A2').Value = '{"Total_reviews":700,"Ease_of_use":4.2,"Features":4,"Design":3.8,"Support":3.9}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":false,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"TeamDynamix","Year_founded":2001,"Country":"United States"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Project Management Integration"]'

# sheet15.xml
$ws = $wb.Worksheets.Item(15)
$ws.Range('A2').Value = '{"Total_reviews":500,"Ease_of_use":4.3,"Features":4.1,"Design":4,"Support":4}'
$ws.Range('B2').Value = '{"API_access":true,"Integrations":{"Active Directory":true,"Answer GPT":false,"Assess360":false,"BigID":false,"Cozyroc SSIS+ Suite":false,"CloudHub":false,"Elastic Observability":false,"Exalate":false,"Incydr":false,"Nexpose":false,"Other_available_integrations":"Limited integrations available"}}'
$ws.Range('C2').Value = '{"Pricing_tiers":"Multiple tiers based on users and features","Free_version_availability":false,"Free_trial_availability":true}'
$ws.Range('D2').Value = '{"SaaS":true,"iPhone":true,"iPad":true,"Android":true,"Windows":true,"Mac":false,"Linux":false}'
$ws.Range('E2').Value = '{"Phone_Support":true,"24/7_Live_Support":false,"Online_Support":true}'
$ws.Range('F2').Value = '{"Documentation":true,"Webinars":true,"Live_online_sessions":true,"In_person_training":false}'
$ws.Range('G2').Value = '{"Company_name":"InvGate","Year_founded":2009,"Country":"Argentina"}'
$ws.Range('H2').Value = '["Incident Management","Problem Management","Change Management","Knowledge Management","Asset Management","Service Catalog","Request Management","IT Service Continuity Management","Workflow Automation","Reporting and Analytics","Self-service Portal","ITIL Process Management"]'

